$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value() = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '48.065.12'
$ws.Range("E2").Value() = '  -0.42%  '
Set-TextValue $ws.Range("D3") '2.494.85'
$ws.Range("E3").Value() = '  -1.11%  '
$ws.Range("E4").Value() = '  -0.05%  '
Set-TextValue $ws.Range("D5") '319.52'
$ws.Range("E5").Value() = '  -1.21%  '
Set-TextValue $ws.Range("D6") '105.73'
$ws.Range("E6").Value() = '  -3.01%  '
$ws.Range("E7").Value() = '  -1.24%  '
Set-TextValue $ws.Range("D9") '0.537'
$ws.Range("E9").Value() = '  -4.21%  '
Set-TextValue $ws.Range("D10") '38.79'
$ws.Range("E10").Value() = '  -4.04%  '
Set-TextValue $ws.Range("D11") '20.02'
$ws.Range("E11").Value() = '  -1.18%  '
$ws.Range("E12").Value() = '  -2.19%  '
$ws.Range("E13").Value() = '  -0.46%  '
Set-TextValue $ws.Range("D14") '7.08'
$ws.Range("E14").Value() = '  -2.64%  '
Set-TextValue $ws.Range("D15") '2.887.87'
$ws.Range("E15").Value() = '  -0.99%  '
Set-TextValue $ws.Range("D16") '2.495.22'
$ws.Range("E16").Value() = '  -0.76%  '
Set-TextValue $ws.Range("D17") '0.830'
$ws.Range("E17").Value() = '  -3.76%  '
Set-TextValue $ws.Range("D18") '47.923.92'
$ws.Range("E18").Value() = '  -0.47%  '
Set-TextValue $ws.Range("D19") '12.98'
$ws.Range("E19").Value() = '  -1.83%  '
$ws.Range("E20").Value() = '  +8.99%  '
Set-TextValue $ws.Range("D21") '6.63'
$ws.Range("E21").Value() = '  -0.33%  '
$ws.Range("E22").Value() = '  -1.57%  '
Set-TextValue $ws.Range("D23") '71.05'
$ws.Range("E23").Value() = '  -1.91%  '
Set-TextValue $ws.Range("D24") '271.05'
$ws.Range("E24").Value() = '  +1.22%  '
Set-TextValue $ws.Range("D25") '2.52'
$ws.Range("E25").Value() = '  -2.21%  '
$ws.Range("E26").Value() = '  -0.03%  '
Set-TextValue $ws.Range("D27") '25.77'
$ws.Range("E27").Value() = '  -1.60%  '
Set-TextValue $ws.Range("D29") '9.72'
$ws.Range("E29").Value() = '  -4.41%  '
$ws.Range("E30").Value() = '  -3.92%  '
Set-TextValue $ws.Range("D31") '34.79'
$ws.Range("E31").Value() = '  -0.53%  '
Set-TextValue $ws.Range("D32") '49.32'
$ws.Range("E33").Value() = '  -0.05%  '
Set-TextValue $ws.Range("D34") '19.04'
$ws.Range("E34").Value() = '  -4.79%  '
Set-TextValue $ws.Range("D35") '5.28'
$ws.Range("E35").Value() = '  -2.13%  '
Set-TextValue $ws.Range("D36") '0.0771'
$ws.Range("E36").Value() = '  -2.72%  '
Set-TextValue $ws.Range("D37") '1.93'
$ws.Range("E37").Value() = '  -2.79%  '
$ws.Range("E38").Value() = '  -3.01%  '
$ws.Range("E39").Value() = '  -4.52%  '
$ws.Range("B40").Value() = 'Monero'
$ws.Range("C40").Value() = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D40") '122.04'
$ws.Range("E40").Value() = '  +2.58%  '
$ws.Range("B41").Value() = 'WEMIXToken'
$ws.Range("C41").Value() = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D41") '2.23'
$ws.Range("E41").Value() = '  +1.55%  '
$ws.Range("E42").Value() = '  -2.23%  '
$ws.Range("B43").Value() = 'EnergySwap'
$ws.Range("C43").Value() = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D43") '22.22'
$ws.Range("E43").Value() = '  +0.37%  '
Set-TextValue $ws.Range("D44") '0.0301'
$ws.Range("E44").Value() = '  +0.49%  '
Set-TextValue $ws.Range("D45") '1.996.21'
$ws.Range("E45").Value() = '  -0.33%  '
Set-TextValue $ws.Range("D46") '3.13'
$ws.Range("E46").Value() = '  +0.03%  '
Set-TextValue $ws.Range("D47") '1.87'
$ws.Range("E47").Value() = '  -1.25%  '
$ws.Range("E48").Value() = '  -1.23%  '
$ws.Range("E49").Value() = '  -2.00%  '
Set-TextValue $ws.Range("D50") '5.15'
$ws.Range("E50").Value() = '  -2.32%  '
Set-TextValue $ws.Range("D51") '78.91'
$ws.Range("E51").Value() = '  -2.08%  '
